# where-to-focus-v3: drop the "Other / Self nominated category" row and
# annotate two Social subcategories with extra description text, then
# autosize the first two columns and set print/page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the last data row (A13:B13 "Other" / "Self nominated category").
#    Deleting it drops the only references to those two shared strings so
#    they fall out of the shared-string table on save, and the remaining
#    "Environmental" entries shift down to fill the gap.
$ws.Rows.Item(13).Delete()

# 2. Add supporting description text for the "Diversity and inclusion" and
#    "Talent, training and career development" subcategories, in column C,
#    with an explicit black font color (this is what seeds the new font /
#    cell style used for these two cells).
$ws.Range("C4").Value = "A workplace which fosters and supports diversity, inclusiveness"
$ws.Range("C4").Font.Color = 0

$ws.Range("C6").Value = "Support and develop employees to help individuals reach their full potential and improve outcomes for the company"
$ws.Range("C6").Font.Color = 0

# 3. Resize columns A and B to fit their (now final) contents.
$ws.Columns.Item(1).ColumnWidth = 14.140625
$ws.Columns.Item(2).ColumnWidth = 37.42578125

# 4. Move the active selection (matches the saved sheet view state).
$ws.Range("F19").Select() | Out-Null

# 5. Configure page setup for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
